# V 0.52-B49 (pre-Release for testing)
# Add Compass degm and arrow (N-up) Item:
# Insert a new "COMPASS" column into the Tabelle2 (sheet2) item matrix,
# right before the END_OF_COL marker column (which currently sits at EE).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle2")

# Inserting a whole column at EE shifts EE..EG (END_OF_COL, Title) one
# column to the right (EF, EG) and shifts all column widths / row spans /
# the used-range dimension accordingly - exactly like a user selecting
# column EE and doing Insert > Entire Column in Excel.
$ws.Columns("EE").Insert()

# Header row: label the new column like the other item headers.
$ws.Range("EE1").Value = "COMPASS"

# Data rows: mark every aircraft row with the same "|" flag used for all
# the other boolean/flag item columns in this sheet.
for ($r = 2; $r -le 40; $r++) {
    $ws.Range("EE$r").Value = "|"
}

# Restore the user's selection as left by the edit session.
$ws.Range("EJ12").Select() | Out-Null
